$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update the "ffa9b05e..." row (row 3) Correspond Handoff
# Datetime (E3) and Correspond Handback DateTime (H3) to reflect the new
# handback run.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-11 14:33:03"
$wsZh.Range("H3").Value = "2016-03-11 14:33:21"

# "de-de" sheet: same update for its "ffa9b05e..." row (row 3).
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-11 14:33:07"
$wsDe.Range("H3").Value = "2016-03-11 14:33:26"
